$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row for 2022-Q4 right after the
#    header, pushing the existing quarters down by one row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()

# Give the new A2 the same style as the (now shifted) data rows, then
# drop the stray formatting Insert() left on B2:D2 so they match the
# plain (unstyled) data cells used elsewhere in the column.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 12
$total.Range("D2").Value = 1.52

# Renumber the 0-based index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

# ---------------------------------------------------------------------
# 2) Add a new "2022-Q4" sheet (holdings detail), positioned right after
#    "总计" and before the existing "2022-Q3" sheet. Duplicating
#    "2022-Q3" gives us an exact style match "for free", then we
#    overwrite the cell values with the 2022-Q4 data.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$ws = $wb.Worksheets.Item(2)
$ws.Name = "2022-Q4"

# Expand the single data row (row 2) down to row 13 - this clones the
# exact formatting of every column into the new rows in one shot.
$ws.Range("A2:H2").Copy($ws.Range("A3:H13"))

# The fund-code (B) and decimal-looking (D:E:F:G) columns are stored as
# text in the source data (e.g. "005668", "0.0120") - force text format
# first so Excel's COM layer doesn't silently coerce them to numbers.
$ws.Range("B2:B13").NumberFormat = "@"
$ws.Range("D2:G13").NumberFormat = "@"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "005668"
$ws.Range("C2").Value = "融通新能源汽车主题精选灵活配置混合A"
$ws.Range("D2").Value = "7.52"
$ws.Range("E2").Value = "94.15"
$ws.Range("F2").Value = "5.49"
$ws.Range("G2").Value = "0.4128"
$ws.Range("H2").Value = 7
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "014647"
$ws.Range("C3").Value = "融通先进制造混合A"
$ws.Range("D3").Value = "4.65"
$ws.Range("E3").Value = "93.63"
$ws.Range("F3").Value = "5.86"
$ws.Range("G3").Value = "0.2725"
$ws.Range("H3").Value = 5
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "001471"
$ws.Range("C4").Value = "融通新能源灵活配置混合A"
$ws.Range("D4").Value = "5.79"
$ws.Range("E4").Value = "91.82"
$ws.Range("F4").Value = "3.43"
$ws.Range("G4").Value = "0.1986"
$ws.Range("H4").Value = 9
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "014648"
$ws.Range("C5").Value = "融通先进制造混合C"
$ws.Range("D5").Value = "3.07"
$ws.Range("E5").Value = "93.63"
$ws.Range("F5").Value = "5.86"
$ws.Range("G5").Value = "0.1799"
$ws.Range("H5").Value = 5
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "001541"
$ws.Range("C6").Value = "汇添富民营新动力股票"
$ws.Range("D6").Value = "3.63"
$ws.Range("E6").Value = "81.10"
$ws.Range("F6").Value = "4.85"
$ws.Range("G6").Value = "0.1761"
$ws.Range("H6").Value = 2
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "660005"
$ws.Range("C7").Value = "农银中小盘混合"
$ws.Range("D7").Value = "6.71"
$ws.Range("E7").Value = "75.58"
$ws.Range("F7").Value = "2.03"
$ws.Range("G7").Value = "0.1362"
$ws.Range("H7").Value = 6
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "003655"
$ws.Range("C8").Value = "信澳新财富灵活配置混合"
$ws.Range("D8").Value = "3.03"
$ws.Range("E8").Value = "78.65"
$ws.Range("F8").Value = "2.24"
$ws.Range("G8").Value = "0.0679"
$ws.Range("H8").Value = 7
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "009835"
$ws.Range("C9").Value = "融通新能源汽车主题精选灵活配置混合C"
$ws.Range("D9").Value = "0.72"
$ws.Range("E9").Value = "94.15"
$ws.Range("F9").Value = "5.49"
$ws.Range("G9").Value = "0.0395"
$ws.Range("H9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "001983"
$ws.Range("C10").Value = "中邮低碳经济灵活配置混合"
$ws.Range("D10").Value = "0.48"
$ws.Range("E10").Value = "91.35"
$ws.Range("F10").Value = "4.29"
$ws.Range("G10").Value = "0.0206"
$ws.Range("H10").Value = 9
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "012005"
$ws.Range("C11").Value = "信澳恒盛混合A"
$ws.Range("D11").Value = "1.29"
$ws.Range("E11").Value = "35.97"
$ws.Range("F11").Value = "0.93"
$ws.Range("G11").Value = "0.0120"
$ws.Range("H11").Value = 10
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "014948"
$ws.Range("C12").Value = "融通新能源灵活配置混合C"
$ws.Range("D12").Value = "0.11"
$ws.Range("E12").Value = "91.82"
$ws.Range("F12").Value = "3.43"
$ws.Range("G12").Value = "0.0038"
$ws.Range("H12").Value = 9
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "012006"
$ws.Range("C13").Value = "信澳恒盛混合C"
$ws.Range("D13").Value = "0.09"
$ws.Range("E13").Value = "35.97"
$ws.Range("F13").Value = "0.93"
$ws.Range("G13").Value = "0.0008"
$ws.Range("H13").Value = 10

# Remove the temporary "text" number format so these cells end up with
# no explicit style, matching the rest of the workbook's data cells.
$ws.Range("B2:B13").ClearFormats()
$ws.Range("D2:G13").ClearFormats()
